$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -21.9607
$ws.Range("A12").Value = -21.52600000000001
$ws.Range("A18").Value = -22.2827
$ws.Range("A37").Value = -19.82819999999999
$ws.Range("A55").Value = -22.54710000000001
$ws.Range("A68").Value = -21.546
$ws.Range("A77").Value = -20.6708
$ws.Range("A78").Value = -20.30439999999998
